$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Append 7 new data rows (268-274) to the "Landscaping Data" table on Sheet1.
# Columns: A=Date B=Plant_Type C=Plant_Size D=Low E=High F=Temp_Diff(formula)
#          G=Rain H=Growth I=Pruned J=Quadrant K=Shade L=UV M=Humidity
#          N=Dew_Point O=Pressure P=Wind_Gust Q=Cloud_Cover R=Visibility
#          S=AQI T=Pollen
# ---------------------------------------------------------------------------

$newRows = @(
    @{ Row=268; A=45825; B="Flowering";    C="Large";  D=70; E=77; G=2.52; H=0.3;                 I="No"; J=2; K="Neutral"; L=8; M=0.84; N=72; O=29.88; P=8; Q=0.76; R=9.9; S=31; T=30 },
    @{ Row=269; A=45825; B="Nonflowering"; C="Medium"; D=70; E=77; G=2.52; H=0.6;                 I="No"; J=3; K="Dark";    L=8; M=0.84; N=72; O=29.88; P=8; Q=0.76; R=9.9; S=31; T=30 },
    @{ Row=270; A=45825; B="Nonflowering"; C="Small";  D=70; E=77; G=2.52; H=0.55000000000000004; I="No"; J=3; K="Neutral"; L=8; M=0.84; N=72; O=29.88; P=8; Q=0.76; R=9.9; S=31; T=30 },
    @{ Row=271; A=45825; B="Nonflowering"; C="Medium"; D=70; E=77; G=2.52; H=0.8;                 I="No"; J=3; K="Dark";    L=8; M=0.84; N=72; O=29.88; P=8; Q=0.76; R=9.9; S=31; T=30 },
    @{ Row=272; A=45825; B="Nonflowering"; C="Medium"; D=70; E=77; G=2.52; H=0.75;                I="No"; J=3; K="Bright";  L=8; M=0.84; N=72; O=29.88; P=8; Q=0.76; R=9.9; S=31; T=30 },
    @{ Row=273; A=45825; B="Nonflowering"; C="Large";  D=70; E=77; G=2.52; H=2;                   I="No"; J=4; K="Neutral"; L=8; M=0.84; N=72; O=29.88; P=8; Q=0.76; R=9.9; S=31; T=30 },
    @{ Row=274; A=45825; B="Tree";         C="Medium"; D=70; E=77; G=2.52; H=5;                   I="No"; J=1; K="Neutral"; L=8; M=0.84; N=72; O=29.88; P=8; Q=0.76; R=9.9; S=31; T=30 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A    # Date
    $ws.Cells.Item($row, 2).Value = $r.B    # Plant_Type
    $ws.Cells.Item($row, 3).Value = $r.C    # Plant_Size
    $ws.Cells.Item($row, 4).Value = $r.D    # Low
    $ws.Cells.Item($row, 5).Value = $r.E    # High
    # F (Temp_Diff) filled in below via the shared formula
    $ws.Cells.Item($row, 7).Value = $r.G    # Rain
    $ws.Cells.Item($row, 8).Value = $r.H    # Growth
    $ws.Cells.Item($row, 9).Value = $r.I    # Pruned
    $ws.Cells.Item($row, 10).Value = $r.J   # Quadrant
    $ws.Cells.Item($row, 11).Value = $r.K   # Shade
    $ws.Cells.Item($row, 12).Value = $r.L   # UV
    $ws.Cells.Item($row, 13).Value = $r.M   # Humidity
    $ws.Cells.Item($row, 14).Value = $r.N   # Dew_Point
    $ws.Cells.Item($row, 15).Value = $r.O   # Pressure
    $ws.Cells.Item($row, 16).Value = $r.P   # Wind_Gust
    $ws.Cells.Item($row, 17).Value = $r.Q   # Cloud_Cover
    $ws.Cells.Item($row, 18).Value = $r.R   # Visibility
    $ws.Cells.Item($row, 19).Value = $r.S   # AQI
    $ws.Cells.Item($row, 20).Value = $r.T   # Pollen
}

# Column A (Date) keeps the existing date number format - copy it from the
# last pre-existing row so the new cells share the same style instead of
# creating a new custom numeric format.
$ws.Range("A267").Copy()
$ws.Range("A268:A274").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Re-apply the values after the paste (paste-formats only touches formats,
# but do this defensively in case values were cleared).
foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
}

# Extend the Temp_Diff formula (=ABS(Low-High)) down through the new rows.
$ws.Range("F268:F274").Formula = "=ABS(D268-E268)"
